$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.836.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.100.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.89"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.17"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0765"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.922"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.406.01"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.147.44"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.811.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.82"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.48"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.87"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.78"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.26"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.75%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.17%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0824"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.98%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0221"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0955"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.74"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.82"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.411.48"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +11.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.17"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.03%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +11.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.47"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.293.28"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.83%  "
